$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 54 values (quarter 01-01-2021) ---
$ws.Range("B54").Value = 126.9
$ws.Range("D54").Value = 133.2
$ws.Range("E54").Value = 113.9
$ws.Range("F54").Value = 107.4
$ws.Range("G54").Value = 113.4
$ws.Range("I54").Value = 93.5
$ws.Range("J54").Value = 117
$ws.Range("L54").Value = 115.5

# --- Append new row 55 (quarter 01-04-2021) ---
# Force the date-like label to be stored as text (matching the rest of
# column A) instead of letting Excel auto-convert it to a date serial.
$ws.Range("A55").NumberFormat = "@"
$ws.Range("A55").Value = "01-04-2021"

# Re-align the cell's style with the rest of the column (no explicit
# style / general format) by copying the format from the cell above.
$ws.Range("A54").Copy()
$ws.Range("A55").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("B55").Value = 130.2
$ws.Range("C55").Value = 129.2
$ws.Range("D55").Value = 129.7
$ws.Range("E55").Value = 117.5
$ws.Range("F55").Value = 106.9
$ws.Range("G55").Value = 107.8
$ws.Range("H55").Value = 88.5
$ws.Range("I55").Value = 112.8
$ws.Range("J55").Value = 122
$ws.Range("K55").Value = 106.4
$ws.Range("L55").Value = 119.5
